$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1727455869996703, 0.2925651826085037, 0.1673036931228319, 0.409027741263147,  0.3847553842407884, 14),
    @(0.2046513626860979, 0.2941678900787714, 0.1637830044104861, 0.4047011297371013, 0.3633995406912939, 13),
    @(0.2207881156328718, 0.3023524304801385, 0.1964366846399542, 0.4432117830563107, 0.4013920676417463, 12),
    @(0.2661302185036374, 0.2968048388866972, 0.2374595161704083, 0.4872981799375084, 0.4281327425770843, 11),
    @(0.2665145217128843, 0.2913645221026481, 0.1729824392181902, 0.4159115762012284, 0.3365716845130267, 10),
    @(0.2914848454201801, 0.2916783576633104, 0.2212395806333997, 0.4703611172635336, 0.3915490853191271, 9),
    @(0.3135314731477879, 0.344869875023549,  0.201523643479535,  0.4489138486163408, 0.343464132168225,  8),
    @(0.2942408256753315, 0.3552296618050602, 0.3108043548782771, 0.557498300336671,  0.5114663298936425, 7),
    @(0.2556306440737632, 0.3520957816061401, 0.2176259007180903, 0.4665039128647158, 0.4274747354338359, 6),
    @(0.3025019693341894, 0.3751403399158096, 0.3407721918770704, 0.5837569630223441, 0.5581943550704445, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
